$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: rotate the three rows 34-36 (F:V only) ---
# new34 = old36, new35 = old34, new36 = old35
$cols = 6..22
$r34 = @{}
$r35 = @{}
$r36 = @{}
foreach ($c in $cols) {
  $r34[$c] = $ws.Cells.Item(34,$c).Value()
  $r35[$c] = $ws.Cells.Item(35,$c).Value()
  $r36[$c] = $ws.Cells.Item(36,$c).Value()
}
foreach ($c in $cols) {
  $ws.Cells.Item(34,$c).Value() = $r36[$c]
  $ws.Cells.Item(35,$c).Value() = $r34[$c]
  $ws.Cells.Item(36,$c).Value() = $r35[$c]
}

# --- Step 2: rotate the three rows 91-93 (F:V only) ---
# new91 = old92, new92 = old93, new93 = old91
$s91 = @{}
$s92 = @{}
$s93 = @{}
foreach ($c in $cols) {
  $s91[$c] = $ws.Cells.Item(91,$c).Value()
  $s92[$c] = $ws.Cells.Item(92,$c).Value()
  $s93[$c] = $ws.Cells.Item(93,$c).Value()
}
foreach ($c in $cols) {
  $ws.Cells.Item(91,$c).Value() = $s92[$c]
  $ws.Cells.Item(92,$c).Value() = $s93[$c]
  $ws.Cells.Item(93,$c).Value() = $s91[$c]
}

# --- Step 3: append new rows 125-133 ---
# sheet row 125 (Indice 124)
$ws.Cells.Item(125,1).Value() = 124
$ws.Cells.Item(125,2).Value() = "poland"
$ws.Cells.Item(125,3).Value() = "division-2"
$ws.Cells.Item(125,4).Value() = "2023-2024"
$ws.Cells.Item(125,5).Value() = 45226.85416666666
$ws.Cells.Item(125,6).Value() = "Skra"
$ws.Cells.Item(125,7).Value() = 1
$ws.Cells.Item(125,8).Value() = "S. Wola"
$ws.Cells.Item(125,9).Value() = 2
$ws.Cells.Item(125,10).Value() = 2.33
$ws.Cells.Item(125,11).Value() = "26/10/2023 08:42"
$ws.Cells.Item(125,12).Value() = 2.16
$ws.Cells.Item(125,13).Value() = "27/10/2023 20:26"
$ws.Cells.Item(125,14).Value() = 2.96
$ws.Cells.Item(125,15).Value() = "26/10/2023 08:42"
$ws.Cells.Item(125,16).Value() = 3.04
$ws.Cells.Item(125,17).Value() = "27/10/2023 20:26"
$ws.Cells.Item(125,18).Value() = 2.88
$ws.Cells.Item(125,19).Value() = "26/10/2023 08:42"
$ws.Cells.Item(125,20).Value() = 3.59
$ws.Cells.Item(125,21).Value() = "27/10/2023 20:26"
$ws.Cells.Item(125,22).Value() = "https://www.betexplorer.com/football/poland/division-2/skra-czestochowa-stal-stalowa-wola/fmczK0y5/"

# sheet row 126 (Indice 125)
$ws.Cells.Item(126,1).Value() = 125
$ws.Cells.Item(126,2).Value() = "poland"
$ws.Cells.Item(126,3).Value() = "division-2"
$ws.Cells.Item(126,4).Value() = "2023-2024"
$ws.Cells.Item(126,5).Value() = 45227.58333333334
$ws.Cells.Item(126,6).Value() = "Olimpia Elblag"
$ws.Cells.Item(126,7).Value() = 4
$ws.Cells.Item(126,8).Value() = "Sandecja Nowy S."
$ws.Cells.Item(126,9).Value() = 1
$ws.Cells.Item(126,10).Value() = 1.94
$ws.Cells.Item(126,11).Value() = "27/10/2023 02:13"
$ws.Cells.Item(126,12).Value() = 2.14
$ws.Cells.Item(126,13).Value() = "28/10/2023 13:21"
$ws.Cells.Item(126,14).Value() = 3.12
$ws.Cells.Item(126,15).Value() = "27/10/2023 02:13"
$ws.Cells.Item(126,16).Value() = 3.38
$ws.Cells.Item(126,17).Value() = "28/10/2023 13:21"
$ws.Cells.Item(126,18).Value() = 3.57
$ws.Cells.Item(126,19).Value() = "27/10/2023 02:13"
$ws.Cells.Item(126,20).Value() = 3.25
$ws.Cells.Item(126,21).Value() = "28/10/2023 13:21"
$ws.Cells.Item(126,22).Value() = "https://www.betexplorer.com/football/poland/division-2/olimpia-elblag-sandecja-nowy-s/SKwRNrxt/"

# sheet row 127 (Indice 126)
$ws.Cells.Item(127,1).Value() = 126
$ws.Cells.Item(127,2).Value() = "poland"
$ws.Cells.Item(127,3).Value() = "division-2"
$ws.Cells.Item(127,4).Value() = "2023-2024"
$ws.Cells.Item(127,5).Value() = 45227.61458333334
$ws.Cells.Item(127,6).Value() = "Kotwica Kolobrzeg"
$ws.Cells.Item(127,7).Value() = 2
$ws.Cells.Item(127,8).Value() = "KKS Kalisz"
$ws.Cells.Item(127,9).Value() = 1
$ws.Cells.Item(127,10).Value() = 2.17
$ws.Cells.Item(127,11).Value() = "27/10/2023 03:12"
$ws.Cells.Item(127,12).Value() = 2.07
$ws.Cells.Item(127,13).Value() = "28/10/2023 14:35"
$ws.Cells.Item(127,14).Value() = 3.15
$ws.Cells.Item(127,15).Value() = "27/10/2023 03:12"
$ws.Cells.Item(127,16).Value() = 3.22
$ws.Cells.Item(127,17).Value() = "28/10/2023 14:35"
$ws.Cells.Item(127,18).Value() = 2.96
$ws.Cells.Item(127,19).Value() = "27/10/2023 03:12"
$ws.Cells.Item(127,20).Value() = 3.61
$ws.Cells.Item(127,21).Value() = "28/10/2023 14:35"
$ws.Cells.Item(127,22).Value() = "https://www.betexplorer.com/football/poland/division-2/kotwica-kolobrzeg-kks-kalisz/zosZLM6h/"

# sheet row 128 (Indice 127)
$ws.Cells.Item(128,1).Value() = 127
$ws.Cells.Item(128,2).Value() = "poland"
$ws.Cells.Item(128,3).Value() = "division-2"
$ws.Cells.Item(128,4).Value() = "2023-2024"
$ws.Cells.Item(128,5).Value() = 45227.70833333334
$ws.Cells.Item(128,6).Value() = "Polonia Bytom"
$ws.Cells.Item(128,7).Value() = 1
$ws.Cells.Item(128,8).Value() = "Chojniczanka"
$ws.Cells.Item(128,9).Value() = 1
$ws.Cells.Item(128,10).Value() = 2.63
$ws.Cells.Item(128,11).Value() = "27/10/2023 05:12"
$ws.Cells.Item(128,12).Value() = 2.88
$ws.Cells.Item(128,13).Value() = "28/10/2023 16:39"
$ws.Cells.Item(128,14).Value() = 3.01
$ws.Cells.Item(128,15).Value() = "27/10/2023 05:12"
$ws.Cells.Item(128,16).Value() = 3.25
$ws.Cells.Item(128,17).Value() = "28/10/2023 16:39"
$ws.Cells.Item(128,18).Value() = 2.55
$ws.Cells.Item(128,19).Value() = "27/10/2023 05:12"
$ws.Cells.Item(128,20).Value() = 2.4
$ws.Cells.Item(128,21).Value() = "28/10/2023 16:39"
$ws.Cells.Item(128,22).Value() = "https://www.betexplorer.com/football/poland/division-2/polonia-bytom-chojniczanka/2cCK2a6U/"

# sheet row 129 (Indice 128)
$ws.Cells.Item(129,1).Value() = 128
$ws.Cells.Item(129,2).Value() = "poland"
$ws.Cells.Item(129,3).Value() = "division-2"
$ws.Cells.Item(129,4).Value() = "2023-2024"
$ws.Cells.Item(129,5).Value() = 45227.80763888889
$ws.Cells.Item(129,6).Value() = "Ol. Grudziadz"
$ws.Cells.Item(129,7).Value() = 3
$ws.Cells.Item(129,8).Value() = "Stomil Olsztyn"
$ws.Cells.Item(129,9).Value() = 0
$ws.Cells.Item(129,10).Value() = 2.2
$ws.Cells.Item(129,11).Value() = "27/10/2023 07:42"
$ws.Cells.Item(129,12).Value() = 1.96
$ws.Cells.Item(129,13).Value() = "28/10/2023 19:18"
$ws.Cells.Item(129,14).Value() = 3.16
$ws.Cells.Item(129,15).Value() = "27/10/2023 07:42"
$ws.Cells.Item(129,16).Value() = 3.48
$ws.Cells.Item(129,17).Value() = "28/10/2023 19:17"
$ws.Cells.Item(129,18).Value() = 2.92
$ws.Cells.Item(129,19).Value() = "27/10/2023 07:42"
$ws.Cells.Item(129,20).Value() = 3.66
$ws.Cells.Item(129,21).Value() = "28/10/2023 19:18"
$ws.Cells.Item(129,22).Value() = "https://www.betexplorer.com/football/poland/division-2/ol-grudziadz-stomil-olsztyn/EytwLtMb/"

# sheet row 130 (Indice 129)
$ws.Cells.Item(130,1).Value() = 129
$ws.Cells.Item(130,2).Value() = "poland"
$ws.Cells.Item(130,3).Value() = "division-2"
$ws.Cells.Item(130,4).Value() = "2023-2024"
$ws.Cells.Item(130,5).Value() = 45228.47916666666
$ws.Cells.Item(130,6).Value() = "GKS Jastrzebie"
$ws.Cells.Item(130,7).Value() = 2
$ws.Cells.Item(130,8).Value() = "Pogon Siedlce"
$ws.Cells.Item(130,9).Value() = 1
$ws.Cells.Item(130,10).Value() = 2.54
$ws.Cells.Item(130,11).Value() = "28/10/2023 00:42"
$ws.Cells.Item(130,12).Value() = 3.73
$ws.Cells.Item(130,13).Value() = "29/10/2023 11:28"
$ws.Cells.Item(130,14).Value() = 3.07
$ws.Cells.Item(130,15).Value() = "28/10/2023 00:42"
$ws.Cells.Item(130,16).Value() = 4.07
$ws.Cells.Item(130,17).Value() = "29/10/2023 11:28"
$ws.Cells.Item(130,18).Value() = 2.54
$ws.Cells.Item(130,19).Value() = "28/10/2023 00:42"
$ws.Cells.Item(130,20).Value() = 1.79
$ws.Cells.Item(130,21).Value() = "29/10/2023 11:28"
$ws.Cells.Item(130,22).Value() = "https://www.betexplorer.com/football/poland/division-2/gks-jastrzebie-pogon-siedlce/t4anHbLN/"

# sheet row 131 (Indice 130)
$ws.Cells.Item(131,1).Value() = 130
$ws.Cells.Item(131,2).Value() = "poland"
$ws.Cells.Item(131,3).Value() = "division-2"
$ws.Cells.Item(131,4).Value() = "2023-2024"
$ws.Cells.Item(131,5).Value() = 45228.54166666666
$ws.Cells.Item(131,6).Value() = "Lech Poznan II"
$ws.Cells.Item(131,7).Value() = 4
$ws.Cells.Item(131,8).Value() = "Hutnik Krakow"
$ws.Cells.Item(131,9).Value() = 3
$ws.Cells.Item(131,10).Value() = 3.39
$ws.Cells.Item(131,11).Value() = "28/10/2023 02:13"
$ws.Cells.Item(131,12).Value() = 3.18
$ws.Cells.Item(131,13).Value() = "29/10/2023 12:41"
$ws.Cells.Item(131,14).Value() = 3.43
$ws.Cells.Item(131,15).Value() = "28/10/2023 02:13"
$ws.Cells.Item(131,16).Value() = 3.43
$ws.Cells.Item(131,17).Value() = "29/10/2023 12:41"
$ws.Cells.Item(131,18).Value() = 1.93
$ws.Cells.Item(131,19).Value() = "28/10/2023 02:13"
$ws.Cells.Item(131,20).Value() = 2.15
$ws.Cells.Item(131,21).Value() = "29/10/2023 12:41"
$ws.Cells.Item(131,22).Value() = "https://www.betexplorer.com/football/poland/division-2/lech-poznan-hutnik-krakow/tOsVM2in/"

# sheet row 132 (Indice 131)
$ws.Cells.Item(132,1).Value() = 131
$ws.Cells.Item(132,2).Value() = "poland"
$ws.Cells.Item(132,3).Value() = "division-2"
$ws.Cells.Item(132,4).Value() = "2023-2024"
$ws.Cells.Item(132,5).Value() = 45228.55902777778
$ws.Cells.Item(132,6).Value() = "Stezyca"
$ws.Cells.Item(132,7).Value() = 0
$ws.Cells.Item(132,8).Value() = "LKS Lodz II"
$ws.Cells.Item(132,9).Value() = 0
$ws.Cells.Item(132,10).Value() = 1.96
$ws.Cells.Item(132,11).Value() = "28/10/2023 02:42"
$ws.Cells.Item(132,12).Value() = 2.01
$ws.Cells.Item(132,13).Value() = "29/10/2023 13:21"
$ws.Cells.Item(132,14).Value() = 3.41
$ws.Cells.Item(132,15).Value() = "28/10/2023 02:42"
$ws.Cells.Item(132,16).Value() = 3.58
$ws.Cells.Item(132,17).Value() = "29/10/2023 13:21"
$ws.Cells.Item(132,18).Value() = 3.2
$ws.Cells.Item(132,19).Value() = "28/10/2023 02:42"
$ws.Cells.Item(132,20).Value() = 3.38
$ws.Cells.Item(132,21).Value() = "29/10/2023 13:18"
$ws.Cells.Item(132,22).Value() = "https://www.betexplorer.com/football/poland/division-2/stezyca-lks-lodz/YwdvJKjB/"

# sheet row 133 (Indice 132)
$ws.Cells.Item(133,1).Value() = 132
$ws.Cells.Item(133,2).Value() = "poland"
$ws.Cells.Item(133,3).Value() = "division-2"
$ws.Cells.Item(133,4).Value() = "2023-2024"
$ws.Cells.Item(133,5).Value() = 45228.65625
$ws.Cells.Item(133,6).Value() = "Wisla Pulawy"
$ws.Cells.Item(133,7).Value() = 3
$ws.Cells.Item(133,8).Value() = "Zaglebie II"
$ws.Cells.Item(133,9).Value() = 3
$ws.Cells.Item(133,10).Value() = 1.7
$ws.Cells.Item(133,11).Value() = "28/10/2023 05:12"
$ws.Cells.Item(133,12).Value() = 1.72
$ws.Cells.Item(133,13).Value() = "29/10/2023 15:38"
$ws.Cells.Item(133,14).Value() = 3.62
$ws.Cells.Item(133,15).Value() = "28/10/2023 05:12"
$ws.Cells.Item(133,16).Value() = 3.81
$ws.Cells.Item(133,17).Value() = "29/10/2023 15:38"
$ws.Cells.Item(133,18).Value() = 3.99
$ws.Cells.Item(133,19).Value() = "28/10/2023 05:12"
$ws.Cells.Item(133,20).Value() = 4.38
$ws.Cells.Item(133,21).Value() = "29/10/2023 15:38"
$ws.Cells.Item(133,22).Value() = "https://www.betexplorer.com/football/poland/division-2/wisla-pulawy-zaglebie/S0erIv6H/"

# --- Step 4: apply formatting to new rows (copy styles from an existing data row) ---
$ws.Range("A2:V2").Copy()
$ws.Range("A125:V133").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Host "done"